# Remove the "Nome Completo" column (column B) from the sheet.
# This deletes the header and all student-name values, shifting the
# remaining data columns (danceability ... liveness) one column to the
# left (C:K -> B:J).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("B").Delete()
